$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 595
$ws.Range("I38").Value = 595
$ws.Range("K38").Value = 1785
$ws.Range("M38").Value = -1413
# Row 40
$ws.Range("H40").Value = 6999.375
$ws.Range("I40").Value = 6999.375
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6999.375
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -6824.375
# Row 74
$ws.Range("H74").Value = 5605.1055
$ws.Range("I74").Value = 5243.625
$ws.Range("K74").Value = 5243.625
$ws.Range("M74").Value = -4307.625
# Row 77
$ws.Range("H77").Value = 5605.1055
$ws.Range("I77").Value = 5243.625
$ws.Range("K77").Value = 26218.125
$ws.Range("M77").Value = -21538.125
# Row 98
$ws.Range("H98").Value = 2187.5667
$ws.Range("J98").Value = 1994.5
$ws.Range("L98").Value = 1994.5
$ws.Range("N98").Value = -4990.5
# Row 122
$ws.Range("H122").Value = 2187.5667
$ws.Range("J122").Value = 1994.5
$ws.Range("L122").Value = 5983.5
$ws.Range("N122").Value = -10883.5
# Row 125
$ws.Range("H125").Value = 1023.3
$ws.Range("I125").Value = 988
$ws.Range("J125").Value = 1038.4286
$ws.Range("K125").Value = 8892
$ws.Range("L125").Value = 9345.857399999999
$ws.Range("M125").Value = -6432
$ws.Range("N125").Value = -14265.8574
# Row 138
$ws.Range("H138").Value = 5832.306
$ws.Range("I138").Value = 3535.625
$ws.Range("J138").Value = 6280.439
$ws.Range("K138").Value = 10606.875
$ws.Range("L138").Value = 18841.317
$ws.Range("M138").Value = -5466.875
$ws.Range("N138").Value = -29121.317

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 36636.566
$ws.Range("I32").Value = 32087.346
$ws.Range("K32").Value = 32087.346
$ws.Range("M32").Value = -31800.346
# Row 63
$ws.Range("H63").Value = 10300.625
$ws.Range("J63").Value = 15800
$ws.Range("L63").Value = 15800
$ws.Range("N63").Value = -17172
# Row 66
$ws.Range("H66").Value = 10300.625
$ws.Range("J66").Value = 15800
$ws.Range("L66").Value = 79000
$ws.Range("N66").Value = -85864
# Row 74
$ws.Range("H74").Value = 2607047
$ws.Range("I74").Value = 3473922.8
$ws.Range("J74").Value = 6420
$ws.Range("K74").Value = 3473922.8
$ws.Range("L74").Value = 6420
$ws.Range("M74").Value = -3473048.8
$ws.Range("N74").Value = -8168
# Row 77
$ws.Range("H77").Value = 2607047
$ws.Range("I77").Value = 3473922.8
$ws.Range("J77").Value = 6420
$ws.Range("K77").Value = 17369614
$ws.Range("L77").Value = 32100
$ws.Range("M77").Value = -17365246
$ws.Range("N77").Value = -40836
# Row 97
$ws.Range("H97").Value = 214.23077
$ws.Range("J97").Value = 400
$ws.Range("L97").Value = 400
$ws.Range("N97").Value = -1392
# Row 132
$ws.Range("H132").Value = 690757.25
$ws.Range("I132").Value = 972952.1
$ws.Range("J132").Value = 13489.5
$ws.Range("K132").Value = 2918856.3
$ws.Range("L132").Value = 40468.5
$ws.Range("M132").Value = -2916326.3
$ws.Range("N132").Value = -45528.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 92878.14
$ws.Range("J35").Value = 92878.14
$ws.Range("L35").Value = 92878.14
$ws.Range("N35").Value = -93498.14
# Row 81
$ws.Range("H81").Value = 33404.95
$ws.Range("J81").Value = 33404.95
$ws.Range("L81").Value = 33404.95
$ws.Range("N81").Value = -35526.95
# Row 82
$ws.Range("H82").Value = 35477.332
$ws.Range("I82").Value = 11108
$ws.Range("J82").Value = 69594.39999999999
$ws.Range("K82").Value = 11108
$ws.Range("L82").Value = 69594.39999999999
$ws.Range("M82").Value = -10725
$ws.Range("N82").Value = -70360.39999999999
# Row 84
$ws.Range("H84").Value = 33404.95
$ws.Range("J84").Value = 33404.95
$ws.Range("L84").Value = 100214.85
$ws.Range("N84").Value = -110822.85
# Row 85
$ws.Range("H85").Value = 35477.332
$ws.Range("I85").Value = 11108
$ws.Range("J85").Value = 69594.39999999999
$ws.Range("K85").Value = 11108
$ws.Range("L85").Value = 69594.39999999999
$ws.Range("M85").Value = -9782
$ws.Range("N85").Value = -72246.39999999999
# Row 86
$ws.Range("H86").Value = 4220
$ws.Range("I86").Value = 3701.6667
$ws.Range("J86").Value = 4997.5
$ws.Range("K86").Value = 3701.6667
$ws.Range("L86").Value = 4997.5
$ws.Range("M86").Value = -2578.6667
$ws.Range("N86").Value = -7243.5
# Row 89
$ws.Range("H89").Value = 4220
$ws.Range("I89").Value = 3701.6667
$ws.Range("J89").Value = 4997.5
$ws.Range("K89").Value = 18508.3335
$ws.Range("L89").Value = 24987.5
$ws.Range("M89").Value = -12892.3335
$ws.Range("N89").Value = -36219.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1688.5
$ws.Range("I22").Value = 985.625
$ws.Range("J22").Value = 4500
$ws.Range("K22").Value = 985.625
$ws.Range("L22").Value = 4500
$ws.Range("M22").Value = -635.625
$ws.Range("N22").Value = -5200
# Row 31
$ws.Range("H31").Value = 1612682.9
$ws.Range("I31").Value = 1729553.2
$ws.Range("J31").Value = 1511056.4
$ws.Range("K31").Value = 1729553.2
$ws.Range("L31").Value = 1511056.4
$ws.Range("M31").Value = -1729258.2
$ws.Range("N31").Value = -1511646.4
# Row 34
$ws.Range("H34").Value = 1612682.9
$ws.Range("I34").Value = 1729553.2
$ws.Range("J34").Value = 1511056.4
$ws.Range("K34").Value = 1729553.2
$ws.Range("L34").Value = 1511056.4
$ws.Range("M34").Value = -1729351.2
$ws.Range("N34").Value = -1511460.4
# Row 58
$ws.Range("H58").Value = 479334.03
$ws.Range("J58").Value = 4003
$ws.Range("L58").Value = 4003
$ws.Range("N58").Value = -4409
# Row 107
$ws.Range("H107").Value = 1006.9167
$ws.Range("I107").Value = 1178.0588
$ws.Range("K107").Value = 1178.0588
$ws.Range("M107").Value = 741.9412
# Row 122
$ws.Range("H122").Value = 1265.375
$ws.Range("I122").Value = 1185
$ws.Range("J122").Value = 1399.3334
$ws.Range("K122").Value = 3555
$ws.Range("L122").Value = 4198.0002
$ws.Range("M122").Value = -1105
$ws.Range("N122").Value = -9098.0002
# Row 132
$ws.Range("H132").Value = 20455.334
$ws.Range("I132").Value = 23269.723
$ws.Range("J132").Value = 12012.167
$ws.Range("K132").Value = 69809.16900000001
$ws.Range("L132").Value = 36036.501
$ws.Range("M132").Value = -67279.16900000001
$ws.Range("N132").Value = -41096.501
# Row 136
$ws.Range("H136").Value = 479334.03
$ws.Range("J136").Value = 4003
$ws.Range("L136").Value = 12009
$ws.Range("N136").Value = -17109
# Row 141
$ws.Range("H141").Value = 580413
$ws.Range("J141").Value = 624086.9399999999
$ws.Range("L141").Value = 624086.9399999999
$ws.Range("N141").Value = -634446.9399999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 104
$ws.Range("H104").Value = 1321.4286
$ws.Range("I104").Value = 1150
$ws.Range("K104").Value = 3450
$ws.Range("M104").Value = -829
# Row 122
$ws.Range("H122").Value = 34202
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 139
$ws.Range("H139").Value = 2179.4
$ws.Range("J139").Value = 4250
$ws.Range("L139").Value = 12750
$ws.Range("N139").Value = -23030
# Row 140
$ws.Range("H140").Value = 2890.7827
$ws.Range("I140").Value = 2654.9443
$ws.Range("K140").Value = 7964.8329
$ws.Range("M140").Value = -2784.8329

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 6714.3125
$ws.Range("I132").Value = 3907.8462
$ws.Range("J132").Value = 18875.666
$ws.Range("K132").Value = 11723.5386
$ws.Range("L132").Value = 56626.99800000001
$ws.Range("M132").Value = -9193.5386
$ws.Range("N132").Value = -61686.99800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3671.739
$ws.Range("I46").Value = 1805.5
$ws.Range("J46").Value = 3849.476
$ws.Range("K46").Value = 1805.5
$ws.Range("L46").Value = 3849.476
$ws.Range("M46").Value = -1617.5
$ws.Range("N46").Value = -4225.476000000001
# Row 93
$ws.Range("H93").Value = 2112.6667
$ws.Range("I93").Value = 1893.1666
$ws.Range("K93").Value = 1893.1666
$ws.Range("M93").Value = -645.1666
# Row 136
$ws.Range("H136").Value = 5833.3335
$ws.Range("I136").Value = 5833.3335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 17500.0005
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -14950.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2622.2222
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 2622.2222
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 122
$ws.Range("H122").Value = 1430.6666
$ws.Range("I122").Value = 1160.2307
$ws.Range("K122").Value = 3480.6921
$ws.Range("M122").Value = -1030.6921
